$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 32258406  # H33: 33333684 -> 32258406
$ws.Cells.Item(33, 9).Value = 38461930  # I33: 40000400 -> 38461930
$ws.Cells.Item(33, 11).Value = 38461930  # K33: 40000400 -> 38461930
$ws.Cells.Item(33, 13).Value = -38461701  # M33: -40000171 -> -38461701
$ws.Cells.Item(64, 8).Value = 3777.3684  # H64: 3793.5 -> 3777.3684
$ws.Cells.Item(64, 10).Value = 3784  # J64: 3836.6667 -> 3784
$ws.Cells.Item(64, 12).Value = 3784  # L64: 3836.6667 -> 3784
$ws.Cells.Item(64, 14).Value = -4280  # N64: -4332.6667 -> -4280
$ws.Cells.Item(67, 8).Value = 3777.3684  # H67: 3793.5 -> 3777.3684
$ws.Cells.Item(67, 10).Value = 3784  # J67: 3836.6667 -> 3784
$ws.Cells.Item(67, 12).Value = 3784  # L67: 3836.6667 -> 3784
$ws.Cells.Item(67, 14).Value = -5500  # N67: -5552.6667 -> -5500
$ws.Cells.Item(74, 8).Value = 3452.2856  # H74: 3998.6667 -> 3452.2856
$ws.Cells.Item(74, 9).Value = 3249.6667  # I74: 4246 -> 3249.6667
$ws.Cells.Item(74, 10).Value = 3533.3333  # J74: 3875 -> 3533.3333
$ws.Cells.Item(74, 11).Value = 3249.6667  # K74: 4246 -> 3249.6667
$ws.Cells.Item(74, 12).Value = 3533.3333  # L74: 3875 -> 3533.3333
$ws.Cells.Item(74, 13).Value = -2313.6667  # M74: -3310 -> -2313.6667
$ws.Cells.Item(74, 14).Value = -5405.3333  # N74: -5747 -> -5405.3333
$ws.Cells.Item(77, 8).Value = 3452.2856  # H77: 3998.6667 -> 3452.2856
$ws.Cells.Item(77, 9).Value = 3249.6667  # I77: 4246 -> 3249.6667
$ws.Cells.Item(77, 10).Value = 3533.3333  # J77: 3875 -> 3533.3333
$ws.Cells.Item(77, 11).Value = 16248.3335  # K77: 21230 -> 16248.3335
$ws.Cells.Item(77, 12).Value = 17666.6665  # L77: 19375 -> 17666.6665
$ws.Cells.Item(77, 13).Value = -11568.3335  # M77: -16550 -> -11568.3335
$ws.Cells.Item(77, 14).Value = -27026.6665  # N77: -28735 -> -27026.6665
$ws.Cells.Item(86, 8).Value = 4170.364  # H86: 3999.4546 -> 4170.364
$ws.Cells.Item(86, 10).Value = 3587.4  # J86: 3399.4 -> 3587.4
$ws.Cells.Item(86, 12).Value = 3587.4  # L86: 3399.4 -> 3587.4
$ws.Cells.Item(86, 14).Value = -5833.4  # N86: -5645.4 -> -5833.4
$ws.Cells.Item(89, 8).Value = 4170.364  # H89: 3999.4546 -> 4170.364
$ws.Cells.Item(89, 10).Value = 3587.4  # J89: 3399.4 -> 3587.4
$ws.Cells.Item(89, 12).Value = 17937  # L89: 16997 -> 17937
$ws.Cells.Item(89, 14).Value = -29169  # N89: -28229 -> -29169
$ws.Cells.Item(98, 8).Value = 5501.5  # H98: 1445.1923 -> 5501.5
$ws.Cells.Item(98, 9).Value = 10000  # I98: 1278.2142 -> 10000
$ws.Cells.Item(98, 10).Value = 4002  # J98: 1640 -> 4002
$ws.Cells.Item(98, 11).Value = 10000  # K98: 1278.2142 -> 10000
$ws.Cells.Item(98, 12).Value = 4002  # L98: 1640 -> 4002
$ws.Cells.Item(98, 13).Value = -8502  # M98: 219.7858000000001 -> -8502
$ws.Cells.Item(98, 14).Value = -6998  # N98: -4636 -> -6998
$ws.Cells.Item(103, 8).Value = 10017151  # H103: 4623765 -> 10017151
$ws.Cells.Item(103, 10).Value = 725  # J103: 812.7273 -> 725
$ws.Cells.Item(103, 12).Value = 2175  # L103: 2438.1819 -> 2175
$ws.Cells.Item(103, 14).Value = -3347  # N103: -3610.1819 -> -3347
$ws.Cells.Item(112, 8).Value = 1199.3024  # H112: 1174.7959 -> 1199.3024
$ws.Cells.Item(112, 10).Value = 1233.4147  # J112: 1203.5106 -> 1233.4147
$ws.Cells.Item(112, 12).Value = 3700.2441  # L112: 3610.5318 -> 3700.2441
$ws.Cells.Item(112, 14).Value = -5916.2441  # N112: -5826.531800000001 -> -5916.2441
$ws.Cells.Item(115, 8).Value = 1323.3334  # H115: 476.25 -> 1323.3334
$ws.Cells.Item(115, 10).Value = 3000  # J115: 450 -> 3000
$ws.Cells.Item(115, 12).Value = 9000  # L115: 1350 -> 9000
$ws.Cells.Item(115, 14).Value = -12134  # N115: -4484 -> -12134
$ws.Cells.Item(122, 8).Value = 5501.5  # H122: 1445.1923 -> 5501.5
$ws.Cells.Item(122, 9).Value = 10000  # I122: 1278.2142 -> 10000
$ws.Cells.Item(122, 10).Value = 4002  # J122: 1640 -> 4002
$ws.Cells.Item(122, 11).Value = 30000  # K122: 3834.6426 -> 30000
$ws.Cells.Item(122, 12).Value = 12006  # L122: 4920 -> 12006
$ws.Cells.Item(122, 13).Value = -27550  # M122: -1384.6426 -> -27550
$ws.Cells.Item(122, 14).Value = -16906  # N122: -9820 -> -16906
$ws.Cells.Item(138, 8).Value = 1998.1837  # H138: 1993.9791 -> 1998.1837
$ws.Cells.Item(138, 10).Value = 2621.92  # J138: 2639.5 -> 2621.92
$ws.Cells.Item(138, 12).Value = 7865.76  # L138: 7918.5 -> 7865.76
$ws.Cells.Item(138, 14).Value = -18145.76  # N138: -18198.5 -> -18145.76
$ws.Cells.Item(140, 8).Value = 71142.86  # H140: 74200 -> 71142.86
$ws.Cells.Item(140, 10).Value = 71142.86  # J140: 74200 -> 71142.86
$ws.Cells.Item(140, 12).Value = 71142.86  # L140: 74200 -> 71142.86
$ws.Cells.Item(140, 14).Value = -81502.86  # N140: -84560 -> -81502.86

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10243.633  # H32: 9496.091 -> 10243.633
$ws.Cells.Item(32, 9).Value = 11000.238  # I32: 10233.128 -> 11000.238
$ws.Cells.Item(32, 10).Value = 5704  # J32: 5166 -> 5704
$ws.Cells.Item(32, 11).Value = 11000.238  # K32: 10233.128 -> 11000.238
$ws.Cells.Item(32, 12).Value = 5704  # L32: 5166 -> 5704
$ws.Cells.Item(32, 13).Value = -10713.238  # M32: -9946.128000000001 -> -10713.238
$ws.Cells.Item(32, 14).Value = -6278  # N32: -5740 -> -6278
$ws.Cells.Item(122, 8).Value = 4158.548  # H122: 4163.3096 -> 4158.548
$ws.Cells.Item(122, 9).Value = 4837.25  # I122: 4843.5 -> 4837.25
$ws.Cells.Item(122, 11).Value = 14511.75  # K122: 14530.5 -> 14511.75
$ws.Cells.Item(122, 13).Value = -12061.75  # M122: -12080.5 -> -12061.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 13890654  # H86: 12822196 -> 13890654
$ws.Cells.Item(86, 9).Value = 1714.24  # I86: 1671.3846 -> 1714.24
$ws.Cells.Item(86, 10).Value = 45456428  # J86: 38463250 -> 45456428
$ws.Cells.Item(86, 11).Value = 1714.24  # K86: 1671.3846 -> 1714.24
$ws.Cells.Item(86, 12).Value = 45456428  # L86: 38463250 -> 45456428
$ws.Cells.Item(86, 13).Value = -591.24  # M86: -548.3846000000001 -> -591.24
$ws.Cells.Item(86, 14).Value = -45458674  # N86: -38465496 -> -45458674
$ws.Cells.Item(89, 8).Value = 13890654  # H89: 12822196 -> 13890654
$ws.Cells.Item(89, 9).Value = 1714.24  # I89: 1671.3846 -> 1714.24
$ws.Cells.Item(89, 10).Value = 45456428  # J89: 38463250 -> 45456428
$ws.Cells.Item(89, 11).Value = 8571.200000000001  # K89: 8356.923000000001 -> 8571.200000000001
$ws.Cells.Item(89, 12).Value = 227282140  # L89: 192316250 -> 227282140
$ws.Cells.Item(89, 13).Value = -2955.200000000001  # M89: -2740.923000000001 -> -2955.200000000001
$ws.Cells.Item(89, 14).Value = -227293372  # N89: -192327482 -> -227293372

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 102  # H7: 93.666664 -> 102
$ws.Cells.Item(7, 9).Value = 102  # I7: 93.666664 -> 102
$ws.Cells.Item(7, 11).Value = 102  # K7: 93.666664 -> 102
$ws.Cells.Item(7, 13).Value = 11  # M7: 19.333336 -> 11
$ws.Cells.Item(16, 8).Value = 1245.4706  # H16: 1266.5 -> 1245.4706
$ws.Cells.Item(16, 9).Value = 1027.2727  # I16: 1039.1 -> 1027.2727
$ws.Cells.Item(16, 11).Value = 1027.2727  # K16: 1039.1 -> 1027.2727
$ws.Cells.Item(16, 13).Value = -740.2727  # M16: -752.0999999999999 -> -740.2727
$ws.Cells.Item(19, 8).Value = 114.27273  # H19: 115.14286 -> 114.27273
$ws.Cells.Item(19, 9).Value = 119.375  # I19: 119.09091 -> 119.375
$ws.Cells.Item(19, 11).Value = 119.375  # K19: 119.09091 -> 119.375
$ws.Cells.Item(19, 13).Value = 50.625  # M19: 50.90909000000001 -> 50.625
$ws.Cells.Item(22, 8).Value = 575.61536  # H22: 648.5 -> 575.61536
$ws.Cells.Item(22, 9).Value = 550  # I22: 540 -> 550
$ws.Cells.Item(22, 10).Value = 587  # J22: 757 -> 587
$ws.Cells.Item(22, 11).Value = 550  # K22: 540 -> 550
$ws.Cells.Item(22, 12).Value = 587  # L22: 757 -> 587
$ws.Cells.Item(22, 13).Value = -200  # M22: -190 -> -200
$ws.Cells.Item(22, 14).Value = -1287  # N22: -1457 -> -1287
$ws.Cells.Item(24, 8).Value = 114.27273  # H24: 115.14286 -> 114.27273
$ws.Cells.Item(24, 9).Value = 119.375  # I24: 119.09091 -> 119.375
$ws.Cells.Item(24, 11).Value = 119.375  # K24: 119.09091 -> 119.375
$ws.Cells.Item(24, 13).Value = 50.625  # M24: 50.90909000000001 -> 50.625
$ws.Cells.Item(55, 8).Value = 40081  # H55: 9936.25 -> 40081
$ws.Cells.Item(55, 9).Value = 0  # I55: 1366.6666 -> 0
$ws.Cells.Item(55, 10).Value = 40081  # J55: 12792.777 -> 40081
$ws.Cells.Item(55, 11).Value = 0  # K55: 1366.6666 -> 0
$ws.Cells.Item(55, 12).Value = 40081  # L55: 12792.777 -> 40081
$ws.Cells.Item(55, 13).ClearContents()  # M55: was -1051.6666
$ws.Cells.Item(55, 14).Value = -40711  # N55: -13422.777 -> -40711
$ws.Cells.Item(99, 8).Value = 1344.1428  # H99: 1443.8334 -> 1344.1428
$ws.Cells.Item(99, 9).Value = 1244.8889  # I99: 1351.5 -> 1244.8889
$ws.Cells.Item(99, 10).Value = 1522.8  # J99: 1628.5 -> 1522.8
$ws.Cells.Item(99, 11).Value = 1244.8889  # K99: 1351.5 -> 1244.8889
$ws.Cells.Item(99, 12).Value = 1522.8  # L99: 1628.5 -> 1522.8
$ws.Cells.Item(99, 13).Value = 253.1111000000001  # M99: 146.5 -> 253.1111000000001
$ws.Cells.Item(99, 14).Value = -4518.8  # N99: -4624.5 -> -4518.8
$ws.Cells.Item(113, 8).Value = 1245.4706  # H113: 1266.5 -> 1245.4706
$ws.Cells.Item(113, 9).Value = 1027.2727  # I113: 1039.1 -> 1027.2727
$ws.Cells.Item(113, 11).Value = 1027.2727  # K113: 1039.1 -> 1027.2727
$ws.Cells.Item(113, 13).Value = 1142.7273  # M113: 1130.9 -> 1142.7273
$ws.Cells.Item(122, 8).Value = 1219.0264  # H122: 1203.3243 -> 1219.0264
$ws.Cells.Item(122, 10).Value = 1058.3572  # J122: 1001.3077 -> 1058.3572
$ws.Cells.Item(122, 12).Value = 3175.0716  # L122: 3003.9231 -> 3175.0716
$ws.Cells.Item(122, 14).Value = -8075.071599999999  # N122: -7903.9231 -> -8075.071599999999
$ws.Cells.Item(126, 8).Value = 1344.1428  # H126: 1443.8334 -> 1344.1428
$ws.Cells.Item(126, 9).Value = 1244.8889  # I126: 1351.5 -> 1244.8889
$ws.Cells.Item(126, 10).Value = 1522.8  # J126: 1628.5 -> 1522.8
$ws.Cells.Item(126, 11).Value = 3734.6667  # K126: 4054.5 -> 3734.6667
$ws.Cells.Item(126, 12).Value = 4568.4  # L126: 4885.5 -> 4568.4
$ws.Cells.Item(126, 13).Value = -1264.6667  # M126: -1584.5 -> -1264.6667
$ws.Cells.Item(126, 14).Value = -9508.4  # N126: -9825.5 -> -9508.4
$ws.Cells.Item(140, 8).Value = 27375.555  # H140: 26506.666 -> 27375.555
$ws.Cells.Item(140, 10).Value = 27375.555  # J140: 26506.666 -> 27375.555
$ws.Cells.Item(140, 12).Value = 27375.555  # L140: 26506.666 -> 27375.555
$ws.Cells.Item(140, 14).Value = -37735.555  # N140: -36866.666 -> -37735.555

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 693.6667  # H17: 550.25 -> 693.6667
$ws.Cells.Item(17, 9).Value = 101  # I17: 100.5 -> 101
$ws.Cells.Item(17, 10).Value = 990  # J17: 1000 -> 990
$ws.Cells.Item(17, 11).Value = 303  # K17: 301.5 -> 303
$ws.Cells.Item(17, 12).Value = 2970  # L17: 3000 -> 2970
$ws.Cells.Item(17, 13).Value = -134  # M17: -132.5 -> -134
$ws.Cells.Item(17, 14).Value = -3308  # N17: -3338 -> -3308
$ws.Cells.Item(21, 8).Value = 603.625  # H21: 180 -> 603.625
$ws.Cells.Item(21, 9).Value = 499.83334  # I21: 180 -> 499.83334
$ws.Cells.Item(21, 10).Value = 915  # J21: 0 -> 915
$ws.Cells.Item(21, 11).Value = 1499.50002  # K21: 540 -> 1499.50002
$ws.Cells.Item(21, 12).Value = 2745  # L21: 0 -> 2745
$ws.Cells.Item(21, 13).Value = -1326.50002  # M21: -367 -> -1326.50002
$ws.Cells.Item(21, 14).Value = -3091  # N21: None -> -3091
$ws.Cells.Item(123, 8).Value = 5100  # H123: 3750 -> 5100
$ws.Cells.Item(123, 10).Value = 6900  # J123: 4875 -> 6900
$ws.Cells.Item(123, 12).Value = 20700  # L123: 14625 -> 20700
$ws.Cells.Item(123, 14).Value = -25600  # N123: -19525 -> -25600

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 30900  # H49: 31400 -> 30900
$ws.Cells.Item(49, 10).Value = 30900  # J49: 31400 -> 30900
$ws.Cells.Item(49, 12).Value = 30900  # L49: 31400 -> 30900
$ws.Cells.Item(49, 14).Value = -31268  # N49: -31768 -> -31268
$ws.Cells.Item(70, 8).Value = 17482.412  # H70: 24362.4 -> 17482.412
$ws.Cells.Item(70, 9).Value = 21931.854  # I70: 37470.223 -> 21931.854
$ws.Cells.Item(70, 10).Value = 4875.6665  # J70: 4700.6665 -> 4875.6665
$ws.Cells.Item(70, 11).Value = 21931.854  # K70: 37470.223 -> 21931.854
$ws.Cells.Item(70, 12).Value = 4875.6665  # L70: 4700.6665 -> 4875.6665
$ws.Cells.Item(70, 13).Value = -21661.854  # M70: -37200.223 -> -21661.854
$ws.Cells.Item(70, 14).Value = -5415.6665  # N70: -5240.6665 -> -5415.6665
$ws.Cells.Item(73, 8).Value = 17482.412  # H73: 24362.4 -> 17482.412
$ws.Cells.Item(73, 9).Value = 21931.854  # I73: 37470.223 -> 21931.854
$ws.Cells.Item(73, 10).Value = 4875.6665  # J73: 4700.6665 -> 4875.6665
$ws.Cells.Item(73, 11).Value = 21931.854  # K73: 37470.223 -> 21931.854
$ws.Cells.Item(73, 12).Value = 4875.6665  # L73: 4700.6665 -> 4875.6665
$ws.Cells.Item(73, 13).Value = -20995.854  # M73: -36534.223 -> -20995.854
$ws.Cells.Item(73, 14).Value = -6747.6665  # N73: -6572.6665 -> -6747.6665
$ws.Cells.Item(102, 8).Value = 2619.3928  # H102: 2237.2 -> 2619.3928
$ws.Cells.Item(102, 9).Value = 3401.611  # I102: 2584.1538 -> 3401.611
$ws.Cells.Item(102, 10).Value = 1211.4  # J102: 1234.8889 -> 1211.4
$ws.Cells.Item(102, 11).Value = 3401.611  # K102: 2584.1538 -> 3401.611
$ws.Cells.Item(102, 12).Value = 1211.4  # L102: 1234.8889 -> 1211.4
$ws.Cells.Item(102, 13).Value = -1779.611  # M102: -962.1538 -> -1779.611
$ws.Cells.Item(102, 14).Value = -4455.4  # N102: -4478.8889 -> -4455.4
$ws.Cells.Item(126, 8).Value = 3954.32  # H126: 3319.2122 -> 3954.32
$ws.Cells.Item(126, 9).Value = 2489.0908  # I126: 2082.4614 -> 2489.0908
$ws.Cells.Item(126, 10).Value = 5105.5713  # J126: 4123.1 -> 5105.5713
$ws.Cells.Item(126, 11).Value = 7467.2724  # K126: 6247.3842 -> 7467.2724
$ws.Cells.Item(126, 12).Value = 15316.7139  # L126: 12369.3 -> 15316.7139
$ws.Cells.Item(126, 13).Value = -4997.2724  # M126: -3777.3842 -> -4997.2724
$ws.Cells.Item(126, 14).Value = -20256.7139  # N126: -17309.3 -> -20256.7139
$ws.Cells.Item(138, 8).Value = 57170.855  # H138: 57182.668 -> 57170.855
$ws.Cells.Item(138, 10).Value = 57170.855  # J138: 57182.668 -> 57170.855
$ws.Cells.Item(138, 12).Value = 57170.855  # L138: 57182.668 -> 57170.855
$ws.Cells.Item(138, 14).Value = -67450.85500000001  # N138: -67462.66800000001 -> -67450.85500000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 11222.5  # H40: 11072.5 -> 11222.5
$ws.Cells.Item(40, 9).Value = 14745  # I40: 18660 -> 14745
$ws.Cells.Item(40, 10).Value = 7700  # J40: 6520 -> 7700
$ws.Cells.Item(40, 11).Value = 14745  # K40: 18660 -> 14745
$ws.Cells.Item(40, 12).Value = 7700  # L40: 6520 -> 7700
$ws.Cells.Item(40, 13).Value = -14609  # M40: -18524 -> -14609
$ws.Cells.Item(40, 14).Value = -7972  # N40: -6792 -> -7972
$ws.Cells.Item(68, 8).Value = 1930.4375  # H68: 1585.2142 -> 1930.4375
$ws.Cells.Item(68, 9).Value = 1900.875  # I68: 1349 -> 1900.875
$ws.Cells.Item(68, 10).Value = 1960  # J68: 1624.5834 -> 1960
$ws.Cells.Item(68, 11).Value = 1900.875  # K68: 1349 -> 1900.875
$ws.Cells.Item(68, 12).Value = 1960  # L68: 1624.5834 -> 1960
$ws.Cells.Item(68, 13).Value = -1151.875  # M68: -600 -> -1151.875
$ws.Cells.Item(68, 14).Value = -3458  # N68: -3122.5834 -> -3458
$ws.Cells.Item(71, 8).Value = 1930.4375  # H71: 1585.2142 -> 1930.4375
$ws.Cells.Item(71, 9).Value = 1900.875  # I71: 1349 -> 1900.875
$ws.Cells.Item(71, 10).Value = 1960  # J71: 1624.5834 -> 1960
$ws.Cells.Item(71, 11).Value = 9504.375  # K71: 6745 -> 9504.375
$ws.Cells.Item(71, 12).Value = 9800  # L71: 8122.916999999999 -> 9800
$ws.Cells.Item(71, 13).Value = -5760.375  # M71: -3001 -> -5760.375
$ws.Cells.Item(71, 14).Value = -17288  # N71: -15610.917 -> -17288
$ws.Cells.Item(139, 8).Value = 40563.09  # H139: 41472.184 -> 40563.09
$ws.Cells.Item(139, 10).Value = 40554.4  # J139: 41554.4 -> 40554.4
$ws.Cells.Item(139, 12).Value = 40554.4  # L139: 41554.4 -> 40554.4
$ws.Cells.Item(139, 14).Value = -50834.4  # N139: -51834.4 -> -50834.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(22, 8).Value = 6500  # H22: 500 -> 6500
$ws.Cells.Item(22, 9).Value = 5000  # I22: 0 -> 5000
$ws.Cells.Item(22, 10).Value = 8000  # J22: 500 -> 8000
$ws.Cells.Item(22, 11).Value = 5000  # K22: 0 -> 5000
$ws.Cells.Item(22, 12).Value = 8000  # L22: 500 -> 8000
$ws.Cells.Item(22, 13).Value = -4707  # M22: None -> -4707
$ws.Cells.Item(22, 14).Value = -8586  # N22: -1086 -> -8586
$ws.Cells.Item(28, 8).Value = 5019  # H28: 3814.8 -> 5019
$ws.Cells.Item(28, 9).Value = 0  # I28: 2008.5 -> 0
$ws.Cells.Item(28, 11).Value = 0  # K28: 2008.5 -> 0
$ws.Cells.Item(28, 13).ClearContents()  # M28: was -1660.5
$ws.Cells.Item(126, 8).Value = 3025.75  # H126: 2387.524 -> 3025.75
$ws.Cells.Item(126, 9).Value = 1972.2858  # I126: 1649.5883 -> 1972.2858
$ws.Cells.Item(126, 10).Value = 10400  # J126: 5523.75 -> 10400
$ws.Cells.Item(126, 11).Value = 5916.857400000001  # K126: 4948.7649 -> 5916.857400000001
$ws.Cells.Item(126, 12).Value = 31200  # L126: 16571.25 -> 31200
$ws.Cells.Item(126, 13).Value = -3446.857400000001  # M126: -2478.7649 -> -3446.857400000001
$ws.Cells.Item(126, 14).Value = -36140  # N126: -21511.25 -> -36140
